# Fixes for bw2data 2.5 dev release
# Two new "id" rows are inserted into the parameter tables on the sheet:
#   - one before the existing "location" row that used to be row 17
#   - one before the existing "location" row that used to be row 30
# All subsequent rows shift down accordingly (dimension grows from G38 to G40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first "id" row, originally before row 17 ---
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "id"
$ws.Range("B17").Value = 1

# --- Insert second "id" row, before what was originally row 30 ---
# (After the first insertion above, that row is now at index 31.)
$ws.Rows.Item(31).Insert()
$ws.Range("A31").Value = "id"
$ws.Range("B31").Value = 2

# --- Update the view state to match the saved selection/scroll position ---
$ws.Range("B17").Select()
